$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Quotation"
$ws.Range("A2").Value = "Quotation - MOESCHETQ21002160"
$ws.Range("A3").Value = "Quotation - SCB000ETQ21000015"

$ws.Range("B1:C3").Clear()

$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
